$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.614.52'
$ws.Range("E2").Value = '  -1.89%  '
$ws.Range("D3").Value = '1.587.76'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.81%  '
$ws.Range("E6").Value = '  -2.42%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  -2.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0616'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.93%  '
$ws.Range("E10").Value = '  -3.83%  '
$ws.Range("E11").Value = '  -1.52%  '
$ws.Range("D12").Value = '1.809.67'
$ws.Range("D13").Value = '1.587.45'
$ws.Range("E13").Value = '  -2.34%  '
$ws.Range("E14").Value = '  -2.86%  '
$ws.Range("E15").Value = '  -4.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.65'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").Value = '26.596.84'
$ws.Range("E17").Value = '  -1.84%  '
$ws.Range("D18").Value = '0.0₃0728'
$ws.Range("E18").Value = '  -2.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '208.46'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.91%  '
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("E21").Value = '  -3.39%  '
$ws.Range("E22").Value = '  -2.77%  '
$ws.Range("E23").Value = '  -2.62%  '
$ws.Range("E24").Value = '  -2.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.81'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.78%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.24'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.63%  '
$ws.Range("E28").Value = '  -3.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.27'
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("E31").Value = '  -1.86%  '
$ws.Range("E32").Value = '  -3.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.687'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +24.10%  '
$ws.Range("E34").Value = '  -3.14%  '
$ws.Range("D35").Value = '1.305.88'
$ws.Range("E35").Value = '  -3.01%  '
$ws.Range("E36").Value = '  -0.64%  '
$ws.Range("E38").Value = '  -3.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.829'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.40%  '
$ws.Range("E40").Value = '  +0.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.790'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.60%  '
$ws.Range("E42").Value = '  +2.56%  '
$ws.Range("E43").Value = '  -2.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.63'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.77%  '
$ws.Range("D45").Value = '1.723.19'
$ws.Range("E45").Value = '  -2.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.51'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.24%  '
$ws.Range("E47").Value = '  -1.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.840'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0505'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.71%  '
$ws.Range("E50").Value = '  -1.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.55'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.84%  '
